$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.239.52"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "1.661.48"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.87"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5220"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06334"
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.08"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07721"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.673.06"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.430"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").Value = "1.890.83"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5476"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("D16").Value = "0.0₅8235"
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.96"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "26.281.00"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.655"
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.63"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.14"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.084"
$ws.Range("E23").Value = "  -3.91%  "
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.69"
$ws.Range("E26").Value = "  -3.06%  "
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.15"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.414"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05944"
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.643"
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.308"
$ws.Range("E33").Value = "  -4.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.630"
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9791"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.418"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.784"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5903"
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.945"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8605"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D43").Value = "1.030.24"
$ws.Range("E43").Value = "  -4.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.74"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "1.804.22"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  +6.76%  "
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.069"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05186"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.473"
$ws.Range("E51").Value = "  -0.12%  "
